$wb = $excel.ActiveWorkbook

# --- Sheet "Memória" (Plan "Memória") ---------------------------------
$mem = $wb.Worksheets.Item("Memória")
$plan7 = $wb.Worksheets.Item("Plan7")

# Apply the same fill/border formatting already used elsewhere in the
# workbook (Plan7!H6) to the new row 21 band, so the existing style is
# reused instead of a brand new one being created.
$plan7.Range("H6").Copy()
$mem.Range("G21:AL21").PasteSpecial(-4122)

# Row 21 - spells out "teste de string\0" one character per cell
$mem.Range("G21").Value = "t"
$mem.Range("H21").Value = "e"
$mem.Range("I21").Value = "s"
$mem.Range("J21").Value = "t"
$mem.Range("K21").Value = "e"
$mem.Range("M21").Value = "d"
$mem.Range("N21").Value = "e"
$mem.Range("P21").Value = "s"
$mem.Range("Q21").Value = "t"
$mem.Range("R21").Value = "r"
$mem.Range("S21").Value = "i"
$mem.Range("T21").Value = "n"
$mem.Range("U21").Value = "g"
$mem.Range("V21").Value = "\0"

# Row 22 - index / ascii-ish table under the "string" sample
$mem.Range("G22").Value = 0
$mem.Range("H22").Value = 1
$mem.Range("I22").Value = 2
$mem.Range("J22").Value = 3
$mem.Range("K22").Value = 4
$mem.Range("L22").Value = 5
$mem.Range("M22").Value = 6
$mem.Range("N22").Value = 7
$mem.Range("O22").Value = 8
$mem.Range("P22").Value = 9
$mem.Range("Q22").Value = 10
$mem.Range("R22").Value = 11
$mem.Range("S22").Value = 12
$mem.Range("T22").Value = 13
$mem.Range("U22").Value = 14
$mem.Range("V22").Value = 15
$mem.Range("W22").Value = 16
$mem.Range("X22").Value = 17
$mem.Range("Y22").Value = 18
$mem.Range("Z22").Value = "..."

# Row 18 - labels for the "char" / 1-byte integer explanation
$mem.Range("AP18").Value = "char"
$mem.Range("AS18").Value = "( inteiro de 1 B )"

# Narrow the helper columns used by the new string demo
$mem.Columns("Q:Y").ColumnWidth = 2.14

$mem.Range("V21").Select()

# --- Sheet "Plan5" ------------------------------------------------------
$plan5 = $wb.Worksheets.Item("Plan5")
$plan5.Range("G8").Value = "href;target;caption"

# Make Plan5 the active sheet/tab (was Plan6 before the edit) and move
# the selection to G10.
$plan5.Activate()
$plan5.Range("G10").Select()
